$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.177.48"
$ws.Range("E2").Value = "  +0.63%  "

$ws.Range("D3").Value = "3.498.24"
$ws.Range("E3").Value = "  +0.04%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.04"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.89"
$ws.Range("E6").Value = "  +3.97%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.131"
$ws.Range("E9").Value = "  -1.61%  "

$ws.Range("E10").Value = "  -2.46%  "

$ws.Range("E11").Value = "  -0.57%  "

$ws.Range("D12").Value = "4.099.08"
$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "31.17"
$ws.Range("E13").Value = "  +10.51%  "

$ws.Range("D15").Value = "67.090.82"
$ws.Range("E15").Value = "  +0.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("E16").Value = "  -1.06%  "

$ws.Range("D17").Value = "3.495.50"
$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.28"
$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.58"
$ws.Range("E19").Value = "  +3.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "392.73"
$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.00"
$ws.Range("E21").Value = "  +1.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.32"
$ws.Range("E22").Value = "  -0.20%  "

$ws.Range("E23").Value = "  +0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.537"
$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.71"
$ws.Range("E25").Value = "  +0.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000122"
$ws.Range("E26").Value = "  -0.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.24"
$ws.Range("E27").Value = "  +0.57%  "

$ws.Range("E28").Value = "  -0.53%  "

$ws.Range("E29").Value = "  -0.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.11"
$ws.Range("E30").Value = "  -2.78%  "

$ws.Range("E31").Value = "  -2.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.64"
$ws.Range("E33").Value = "  -1.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.38"
$ws.Range("E34").Value = "  +0.21%  "

$ws.Range("E35").Value = "  +1.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.05"
$ws.Range("E36").Value = "  -0.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.878"
$ws.Range("E37").Value = "  -1.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.94"
$ws.Range("E38").Value = "  +1.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.03"
$ws.Range("E39").Value = "  +2.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.67"
$ws.Range("E40").Value = "  -1.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.35"
$ws.Range("E41").Value = "  +1.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0731"
$ws.Range("E42").Value = "  -1.47%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.807.11"
$ws.Range("E43").Value = "  -0.48%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.07"
$ws.Range("E44").Value = "  -0.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.52"
$ws.Range("E45").Value = "  -0.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.55"
$ws.Range("E46").Value = "  -2.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0301"
$ws.Range("E47").Value = "  -3.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "338.68"
$ws.Range("E48").Value = "  -0.88%  "

$ws.Range("E49").Value = "  -2.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.86"
$ws.Range("E50").Value = "  +0.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.847"
$ws.Range("E51").Value = "  -0.46%  "

